# Weekly update for "Pepino dulce" - Mercado Mayorista Lo Valledor de Santiago.
# A new reporting date (2023-04-30, serial 45041) is inserted as the newest
# block of 4 rows (Especial/Primera/Segunda/Tercera) at the top of the
# historical tail, pushing the existing rows 334:355 down to 338:359.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before the current row 334, shifting rows 334:355
# (and everything below) down to 338:359.
$ws.Rows.Item(334).Resize(4).Insert()

# Shared/boilerplate values for every data row in this block.
$mercadoId = 6
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$categoriaId = 100112043
$categoria = "Pepino dulce"
$variedad = "Cultivar IV Región"
$unidad = "`$/bandeja 18 kilos"
$origen = "Provincia de Limarí"
$kgUnidades = 18
$clasificacion = "Hortaliza"
$fechaNueva = 45041

# New rows: [row, calidad, volumen, precioMinimo, precioMaximo, precioPromedio, precioKg]
$nuevasFilas = @(
    @(334, "Especial", 260, 14000, 14000, 14000, 778),
    @(335, "Primera",  420, 12000, 12000, 12000, 667),
    @(336, "Segunda",  370, 9000,  9000,  9000,  500),
    @(337, "Tercera",  210, 7000,  7000,  7000,  389)
)

foreach ($fila in $nuevasFilas) {
    $r = $fila[0]
    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fechaNueva
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $categoriaId
    $ws.Cells.Item($r, 7).Value = $categoria
    $ws.Cells.Item($r, 8).Value = $variedad
    $ws.Cells.Item($r, 9).Value = $fila[1]
    $ws.Cells.Item($r, 10).Value = $fila[2]
    $ws.Cells.Item($r, 11).Value = $fila[3]
    $ws.Cells.Item($r, 12).Value = $fila[4]
    $ws.Cells.Item($r, 13).Value = $fila[5]
    $ws.Cells.Item($r, 14).Value = $unidad
    $ws.Cells.Item($r, 15).Value = $origen
    $ws.Cells.Item($r, 16).Value = $fila[6]
    $ws.Cells.Item($r, 17).Value = $kgUnidades
    $ws.Cells.Item($r, 18).Value = $clasificacion
}
